$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 value updates.
# Column A: "MuSCs" -> "ECs" (sending cluster)
$ws.Range("A2").Value = "ECs"
# Columns B/C/D keep their displayed text (Efna3 / Epha5 / MuSCs), only
# numeric columns E-T change as per the new TPM-derived calculation.
$ws.Range("B2").Value = "Efna3"
$ws.Range("C2").Value = "Epha5"
$ws.Range("D2").Value = "MuSCs"

$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.02551366666666667
$ws.Range("H2").Value = 0.076541
$ws.Range("M2").Value = 0.0237255
$ws.Range("N2").Value = 0.047451
$ws.Range("Q2").Value = 0.0006053244985
$ws.Range("R2").Value = 0.003631946991

$wb.Save()
